# Generate Report for Handback
# Updates the localization-status workbook so that the "12d552c0..." row
# reflects a failed handback transform, and records the error detail
# message (with the corresponding "Error Detail" column widened) on both
# the zh-cn and de-de language sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhSheet = $wb.Worksheets.Item("zh-cn")
$deSheet = $wb.Worksheets.Item("de-de")

$newStatus = "Handback transform failed"
$zhError = "Handback file name: 4m3vceyq.euv is different with handoff file name: 12d552c0-7f90-4578-9db8-c46e9c1b93f2.220150210f1ef4f4adfd3694e8bd9d8ad49087ae.zh-cn."
$deError = "Handback file name: 4m3vceyq.euv is different with handoff file name: 12d552c0-7f90-4578-9db8-c46e9c1b93f2.220150210f1ef4f4adfd3694e8bd9d8ad49087ae.de-de."

# The 12d552c0-... file's status moves from "Ready for handoff" to
# "Handback transform failed" everywhere it is reported: the per-language
# detail sheets and the roll-up Overview sheet.
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# zh-cn sheet: row 3 is the 12d552c0-... file
$zhSheet.Range("C3").Value = $newStatus
$zhSheet.Range("P3").Value = $zhError
# Widen the "Error Detail" column (P / 16) to fit the new message - match
# the width already used by column A (40 characters).
$zhSheet.Columns.Item(16).ColumnWidth = $zhSheet.Columns.Item(1).ColumnWidth

# de-de sheet: row 3 is the 12d552c0-... file
$deSheet.Range("C3").Value = $newStatus
$deSheet.Range("P3").Value = $deError
$deSheet.Columns.Item(16).ColumnWidth = $deSheet.Columns.Item(1).ColumnWidth
